$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 9).Value = 34261648.29981892
$ws.Cells.Item(4, 10).Value = 33538582.93039495
$ws.Cells.Item(4, 11).Value = 35132281.77107649
$ws.Cells.Item(4, 12).Value = 37391418.68131825
$ws.Cells.Item(4, 13).Value = 39611701.85578573
$ws.Cells.Item(4, 14).Value = 44186559.19346254
$ws.Cells.Item(4, 20).Value = 6764.392556726341
$ws.Cells.Item(4, 21).Value = 6149.918938368929
$ws.Cells.Item(4, 22).Value = 5885.790211270981
$ws.Cells.Item(4, 23).Value = 5700.78040575061
$ws.Cells.Item(4, 24).Value = 5828.678907561173
$ws.Cells.Item(4, 25).Value = 6291.692893843448
$ws.Cells.Item(4, 26).Value = -614.4736183574123
$ws.Cells.Item(4, 27).Value = -878.6023454553606
$ws.Cells.Item(4, 28).Value = -1063.612150975731
$ws.Cells.Item(4, 29).Value = -935.7136491651681
$ws.Cells.Item(4, 30).Value = -472.6996628828929
$ws.Cells.Item(4, 31).Value = -0.09083943801374972
$ws.Cells.Item(4, 32).Value = -0.1298863627572443
$ws.Cells.Item(4, 33).Value = -0.1572369051701623
$ws.Cells.Item(4, 34).Value = -0.138329294362835
$ws.Cells.Item(4, 35).Value = -0.06988057817739335
$ws.Cells.Item(5, 9).Value = 884935348.778636
$ws.Cells.Item(5, 10).Value = 893593475.3384469
$ws.Cells.Item(5, 11).Value = 910245273.6143174
$ws.Cells.Item(5, 12).Value = 925453099.246821
$ws.Cells.Item(5, 13).Value = 940283932.2425988
$ws.Cells.Item(5, 14).Value = 993506882.9244764
$ws.Cells.Item(5, 18).Value = 0.9272727272727272
$ws.Cells.Item(5, 20).Value = 7482.458051025095
$ws.Cells.Item(5, 21).Value = 7129.327514557919
$ws.Cells.Item(5, 22).Value = 6884.714180688796
$ws.Cells.Item(5, 23).Value = 6636.760402936107
$ws.Cells.Item(5, 24).Value = 6490.132677907763
$ws.Cells.Item(5, 25).Value = 6675.761696273266
$ws.Cells.Item(5, 26).Value = -353.1305364671762
$ws.Cells.Item(5, 27).Value = -597.7438703362996
$ws.Cells.Item(5, 28).Value = -845.697648088988
$ws.Cells.Item(5, 29).Value = -992.3253731173327
$ws.Cells.Item(5, 30).Value = -806.6963547518289
$ws.Cells.Item(5, 31).Value = -0.04719445589391547
$ws.Cells.Item(5, 32).Value = -0.07988603026707364
$ws.Cells.Item(5, 33).Value = -0.1130240413406832
$ws.Cells.Item(5, 34).Value = -0.1326202387437888
$ws.Cells.Item(5, 35).Value = -0.1078116775597976
$ws.Cells.Item(6, 9).Value = 226550911.3650893
$ws.Cells.Item(6, 10).Value = 216861398.9467449
$ws.Cells.Item(6, 11).Value = 221391219.8751115
$ws.Cells.Item(6, 12).Value = 224453585.9165613
$ws.Cells.Item(6, 13).Value = 216845196.498805
$ws.Cells.Item(6, 14).Value = 208245257.1827864
$ws.Cells.Item(6, 17).Value = 0.7611464968152867
$ws.Cells.Item(6, 19).Value = 0.643312101910828
$ws.Cells.Item(6, 20).Value = 11516.12206710328
$ws.Cells.Item(6, 21).Value = 10728.80814063944
$ws.Cells.Item(6, 22).Value = 10024.73317825223
$ws.Cells.Item(6, 23).Value = 10012.42716255431
$ws.Cells.Item(6, 24).Value = 9500.753439309719
$ws.Cells.Item(6, 25).Value = 10137.53564320837
$ws.Cells.Item(6, 26).Value = -787.313926463843
$ws.Cells.Item(6, 27).Value = -1491.388888851045
$ws.Cells.Item(6, 28).Value = -1503.694904548965
$ws.Cells.Item(6, 29).Value = -2015.368627793559
$ws.Cells.Item(6, 30).Value = -1378.586423894907
$ws.Cells.Item(6, 31).Value = -0.06836623664426655
$ws.Cells.Item(6, 32).Value = -0.1295044356217199
$ws.Cells.Item(6, 33).Value = -0.1305730258664407
$ws.Cells.Item(6, 34).Value = -0.1750041043373985
$ws.Cells.Item(6, 35).Value = -0.1197092576704227
$ws.Cells.Item(7, 2).Value = 13902
$ws.Cells.Item(7, 3).Value = 3775452.5
$ws.Cells.Item(7, 4).Value = 3842055.5
$ws.Cells.Item(7, 5).Value = 3888311.5
$ws.Cells.Item(7, 6).Value = 3901997
$ws.Cells.Item(7, 7).Value = 3893404.5
$ws.Cells.Item(7, 8).Value = 3849271.5
$ws.Cells.Item(7, 9).Value = 22036363958.62372
$ws.Cells.Item(7, 10).Value = 21564035051.56988
$ws.Cells.Item(7, 11).Value = 21652543143.27745
$ws.Cells.Item(7, 12).Value = 21737816174.22038
$ws.Cells.Item(7, 13).Value = 21275819964.90931
$ws.Cells.Item(7, 14).Value = 21878692778.50014
$ws.Cells.Item(7, 15).Value = 0.7370162566537188
$ws.Cells.Item(7, 16).Value = 0.7394619479211624
$ws.Cells.Item(7, 17).Value = 0.7097539922313336
$ws.Cells.Item(7, 18).Value = 0.7666522802474464
$ws.Cells.Item(7, 19).Value = 0.5880448856279672
$ws.Cells.Item(7, 20).Value = 5836.747769604762
$ws.Cells.Item(7, 21).Value = 5612.629763305054
$ws.Cells.Item(7, 22).Value = 5568.623589770893
$ws.Cells.Item(7, 23).Value = 5570.946408780011
$ws.Cells.Item(7, 24).Value = 5464.579897852718
$ws.Cells.Item(7, 25).Value = 5683.85284812987
$ws.Cells.Item(7, 26).Value = -224.1180062997073
$ws.Cells.Item(7, 27).Value = -268.1241798338688
$ws.Cells.Item(7, 28).Value = -265.8013608247511
$ws.Cells.Item(7, 29).Value = -372.1678717520435
$ws.Cells.Item(7, 30).Value = -152.8949214748918
$ws.Cells.Item(7, 31).Value = -0.03839775421970715
$ws.Cells.Item(7, 32).Value = -0.0459372565712266
$ws.Cells.Item(7, 33).Value = -0.0455392919681965
$ws.Cells.Item(7, 34).Value = -0.06376288413388898
$ws.Cells.Item(7, 35).Value = -0.02619522506542116
$ws.Cells.Item(8, 2).Value = 2561
$ws.Cells.Item(8, 3).Value = 2505033.5
$ws.Cells.Item(8, 4).Value = 2527966.5
$ws.Cells.Item(8, 5).Value = 2551769
$ws.Cells.Item(8, 6).Value = 2600939
$ws.Cells.Item(8, 7).Value = 2656527.5
$ws.Cells.Item(8, 8).Value = 2712370.5
$ws.Cells.Item(8, 9).Value = 17483655312.56708
$ws.Cells.Item(8, 10).Value = 17069797462.35013
$ws.Cells.Item(8, 11).Value = 16726241348.55186
$ws.Cells.Item(8, 12).Value = 16841549108.31321
$ws.Cells.Item(8, 13).Value = 16723713678.39743
$ws.Cells.Item(8, 14).Value = 17898312753.82991
$ws.Cells.Item(8, 15).Value = 0.7270597422881687
$ws.Cells.Item(8, 16).Value = 0.8367825068332683
$ws.Cells.Item(8, 17).Value = 0.8492776259273721
$ws.Cells.Item(8, 18).Value = 0.9140960562280359
$ws.Cells.Item(8, 19).Value = 0.7461928934010152
$ws.Cells.Item(8, 20).Value = 6979.409781373015
$ws.Cells.Item(8, 21).Value = 6752.382779736254
$ws.Cells.Item(8, 22).Value = 6554.763126502384
$ws.Cells.Item(8, 23).Value = 6475.180351524278
$ws.Cells.Item(8, 24).Value = 6295.32864929779
$ws.Cells.Item(8, 25).Value = 6598.771352892205
$ws.Cells.Item(8, 26).Value = -227.0270016367613
$ws.Cells.Item(8, 27).Value = -424.6466548706312
$ws.Cells.Item(8, 28).Value = -504.2294298487368
$ws.Cells.Item(8, 29).Value = -684.0811320752255
$ws.Cells.Item(8, 30).Value = -380.63842848081
$ws.Cells.Item(8, 31).Value = -0.03252810893016511
$ws.Cells.Item(8, 32).Value = -0.06084277441395525
$ws.Cells.Item(8, 33).Value = -0.0722452822865407
$ws.Cells.Item(8, 34).Value = -0.09801418078372959
$ws.Cells.Item(8, 35).Value = -0.05453733774117642
$ws.Cells.Item(9, 9).Value = 1985532039.197385
$ws.Cells.Item(9, 10).Value = 1998972787.702092
$ws.Cells.Item(9, 11).Value = 2072372855.365217
$ws.Cells.Item(9, 12).Value = 2146450563.122761
$ws.Cells.Item(9, 13).Value = 2099815980.546012
$ws.Cells.Item(9, 14).Value = 2292087499.870823
$ws.Cells.Item(9, 15).Value = 0.7223042836041359
$ws.Cells.Item(9, 16).Value = 0.7385524372230429
$ws.Cells.Item(9, 17).Value = 0.7710487444608567
$ws.Cells.Item(9, 18).Value = 0.844903988183161
$ws.Cells.Item(9, 19).Value = 0.7651403249630724
$ws.Cells.Item(9, 20).Value = 26312.90098792562
$ws.Cells.Item(9, 21).Value = 25380.88075905093
$ws.Cells.Item(9, 22).Value = 25126.5236623954
$ws.Cells.Item(9, 23).Value = 24753.93186743121
$ws.Cells.Item(9, 24).Value = 23085.68265998968
$ws.Cells.Item(9, 25).Value = 24135.11250903793
$ws.Cells.Item(9, 26).Value = -932.0202288746877
$ws.Cells.Item(9, 27).Value = -1186.377325530219
$ws.Cells.Item(9, 28).Value = -1558.969120494414
$ws.Cells.Item(9, 29).Value = -3227.218327935934
$ws.Cells.Item(9, 30).Value = -2177.788478887687
$ws.Cells.Item(9, 31).Value = -0.0354206565555949
$ws.Cells.Item(9, 32).Value = -0.04508728726166
$ws.Cells.Item(9, 33).Value = -0.0592473297113757
$ws.Cells.Item(9, 34).Value = -0.1226477585811169
$ws.Cells.Item(9, 35).Value = -0.08276504669276197
